$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.121.68'
$ws.Range("E2").Value = '  -1.74%  '

$ws.Range("D3").Value = '2.340.35'
$ws.Range("E3").Value = '  +0.52%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.81'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.50%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.567'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("E9").Value = '  -5.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.67'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.59%  '

$ws.Range("E11").Value = '  -2.59%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.11'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.77%  '

$ws.Range("E13").Value = '  -1.44%  '

$ws.Range("D14").Value = '2.696.77'
$ws.Range("E14").Value = '  +0.46%  '

$ws.Range("D15").Value = '2.340.28'
$ws.Range("E15").Value = '  +0.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.69'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.10%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.803'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.33%  '

$ws.Range("D18").Value = '46.095.10'
$ws.Range("E18").Value = '  -1.58%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.49%  '

$ws.Range("D20").Value = '0.0₃0962'
$ws.Range("E20").Value = '  +0.65%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.46%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.68'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.94'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.79%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.40%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.91'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.93%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '40.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.43%  '

$ws.Range("E28").Value = '  -2.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.66'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.95%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.93'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.94%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.61'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +15.17%  '

$ws.Range("E32").Value = '  +5.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.48'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.55%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '144.72'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.15%  '

$ws.Range("E35").Value = '  -5.48%  '

$ws.Range("E36").Value = '  -2.02%  '

$ws.Range("E37").Value = '  -3.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.82'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.07%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.11'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.34%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.88'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.27%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0299'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.82%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.86%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.04%  '

$ws.Range("D44").Value = '1.845.43'
$ws.Range("E44").Value = '  +2.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.82'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.82'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.186'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.65%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '69.63'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.15%  '

$ws.Range("D49").Value = '2.569.57'
$ws.Range("E49").Value = '  +0.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '96.41'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.51%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.74'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.02%  '

